# Mini Project 1 Completed
# Update the "Signup" sheet test data: row 2's first name "Hate" -> "Hope"
# and its email "hatebrotest@gmail.com" -> "HopenNopeh@gmail.com".
# Also apply the Hyperlink style to C3 (email column already has a live
# hyperlink defined, just wasn't styled as one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Signup")

$ws.Range("A2").Value = "Hope"
$ws.Range("C2").Value = "HopenNopeh@gmail.com"

$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
